$wb = $excel.ActiveWorkbook

# Sheet "Metadata": update URL, Title, Date, Description values
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "http://fhir.nmdp.org/ig/matchsync/ValueSet/ms-rh-codes"
$wsMeta.Range("B5").Value = "MatchSync Rh Value Sets"
$wsMeta.Range("B8").Value = "2023-01-12T09:36:27-06:00"
$wsMeta.Range("B11").Value = "MatchSync codes for blood Rh. Combines NMDP and LOINC code"

# Sheet "Include ValueSets": update the NMDP ValueSet URL
$wsInclude = $wb.Worksheets.Item("Include ValueSets")
$wsInclude.Range("A2").Value = "http://fhir.nmdp.org/ig/matchsync/ValueSet/nmdp-rh-status-codes"
